# Updates cryptos list: Price (D) and Volume(1h) (E) columns for Sheet1.
# Values are text (not numbers), so each cell is written via a temporary
# text formula + copy/paste-values, which avoids Excel re-typing numeric-
# looking strings (e.g. "240.07", "0.0775") as floating point numbers and
# avoids any NumberFormat/style churn on the target cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Formula = '="44.032.50"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E2")
$c.Formula = '="  +0.11%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D3")
$c.Formula = '="2.360.30"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E3")
$c.Formula = '="  +0.13%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E4")
$c.Formula = '="  +0.09%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E5")
$c.Formula = '="  +0.76%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D6")
$c.Formula = '="240.07"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E6")
$c.Formula = '="  +0.78%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D7")
$c.Formula = '="74.44"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E8")
$c.Formula = '="  +0.07%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E9")
$c.Formula = '="  +11.46%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E10")
$c.Formula = '="  +0.48%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D11")
$c.Formula = '="57.25"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E11")
$c.Formula = '="  -0.07%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D12")
$c.Formula = '="32.31"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E12")
$c.Formula = '="  +10.45%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D13")
$c.Formula = '="7.33"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E13")
$c.Formula = '="  +9.80%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E14")
$c.Formula = '="  +0.79%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D15")
$c.Formula = '="2.711.77"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E15")
$c.Formula = '="  +0.20%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D16")
$c.Formula = '="16.63"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E16")
$c.Formula = '="  -1.15%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E17")
$c.Formula = '="  -0.33%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D18")
$c.Formula = '="2.354.38"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E18")
$c.Formula = '="  -0.56%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D19")
$c.Formula = '="43.940.58"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E19")
$c.Formula = '="  -0.06%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E20")
$c.Formula = '="  -0.32%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E21")
$c.Formula = '="  +5.16%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D22")
$c.Formula = '="77.07"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E22")
$c.Formula = '="  -1.18%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D23")
$c.Formula = '="258.90"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E23")
$c.Formula = '="  +1.56%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E24")
$c.Formula = '="  +24.93%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E26")
$c.Formula = '="  -0.52%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D27")
$c.Formula = '="3.67"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E27")
$c.Formula = '="  -2.12%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D28")
$c.Formula = '="10.82"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E28")
$c.Formula = '="  +3.36%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E29")
$c.Formula = '="  -0.29%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E30")
$c.Formula = '="  +1.64%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D31")
$c.Formula = '="175.82"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E31")
$c.Formula = '="  +1.91%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E32")
$c.Formula = '="  -2.16%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E33")
$c.Formula = '="  +3.45%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D34")
$c.Formula = '="0.0775"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E34")
$c.Formula = '="  +7.40%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E35")
$c.Formula = '="  +1.54%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D36")
$c.Formula = '="5.47"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E36")
$c.Formula = '="  +4.07%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E37")
$c.Formula = '="  -3.81%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E38")
$c.Formula = '="  -3.00%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E39")
$c.Formula = '="  -1.52%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E40")
$c.Formula = '="  +4.97%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E41")
$c.Formula = '="  +15.24%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D42")
$c.Formula = '="0.209"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E42")
$c.Formula = '="  +15.09%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E43")
$c.Formula = '="  +3.48%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E44")
$c.Formula = '="  -1.55%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E45")
$c.Formula = '="  -0.01%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E46")
$c.Formula = '="  +6.82%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D47")
$c.Formula = '="2.55"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E47")
$c.Formula = '="  +9.63%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D48")
$c.Formula = '="58.14"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E48")
$c.Formula = '="  +10.79%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E49")
$c.Formula = '="  -0.02%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E50")
$c.Formula = '="  +0.83%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D51")
$c.Formula = '="100.41"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E51")
$c.Formula = '="  +2.08%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$excel.CutCopyMode = $false
